$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 43 (Order_ID 42): the deleted item's order total/flag change
$ws.Cells.Item(43, 5).Value = 4574.394932735427
$ws.Cells.Item(43, 6).Value = 1

# New row 44 (Order_ID 43)
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 2
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = "2025-03-03 18:27:18"
$ws.Cells.Item(44, 5).Value = 22.522
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = "werf"

# New row 45 (Order_ID 44)
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 2
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = "2025-03-03 18:29:00"
$ws.Cells.Item(45, 5).Value = -8668.918703389829
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = "dfsggs"

# New row 46 (Order_ID 45)
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = "2025-03-03 18:31:04"
$ws.Cells.Item(46, 5).Value = 3342.896860986547
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = "sfdsgdfg"
